$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.139.24'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '1.875.81'
$ws.Range('E3').Value = '  -1.46%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.81'
$ws.Range('E5').Value = '  +0.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9988'
$ws.Range('E6').Value = '  -0.23%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5047'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3846'
$ws.Range('E8').Value = '  -2.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.09022'
$ws.Range('E9').Value = '  -5.83%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.122'
$ws.Range('E10').Value = '  -1.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.72'
$ws.Range('E11').Value = '  -0.77%  '
$ws.Range('E12').Value = '  -0.17%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.79'
$ws.Range('E13').Value = '  -0.33%  '
$ws.Range('D14').Value = '1.876.52'
$ws.Range('E14').Value = '  -0.73%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.276'
$ws.Range('E15').Value = '  -0.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9990'
$ws.Range('E16').Value = '  -0.19%  '
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.32'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06650'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.23'
$ws.Range('E20').Value = '  +1.98%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9994'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.142'
$ws.Range('E22').Value = '  -0.99%  '
$ws.Range('D23').Value = '28.162.19'
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.48'
$ws.Range('E24').Value = '  +1.77%  '
$ws.Range('E25').Value = '  -2.37%  '
$ws.Range('B26').Value = 'LidoDAOToken'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.547'
$ws.Range('E26').Value = '  -4.25%  '
$ws.Range('B27').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C27').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D27').Value = '2.085.93'
$ws.Range('E27').Value = '  -1.05%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.85'
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '157.02'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '126.94'
$ws.Range('E30').Value = '  +0.28%  '
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.064'
$ws.Range('E32').Value = '  -1.82%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.616'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('E34').Value = '  -0.65%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '9.450'
$ws.Range('E35').Value = '  -1.28%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06591'
$ws.Range('E36').Value = '  -0.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02409'
$ws.Range('E37').Value = '  -0.68%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2196'
$ws.Range('E38').Value = '  +0.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.295'
$ws.Range('E39').Value = '  +1.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.214'
$ws.Range('E40').Value = '  -1.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6407'
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '11.52'
$ws.Range('E42').Value = '  +1.47%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.931'
$ws.Range('E43').Value = '  -1.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9987'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.28'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6036'
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.277'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('E48').Value = '  -1.64%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.240'
$ws.Range('E49').Value = '  +4.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.007'
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '121.38'
$ws.Range('E51').Value = '  -1.32%  '
